# Build a "Result" summary sheet (ISBN / Tittel / Forlag) backed by a
# native Excel Table ("MyTable"), placed after the existing Sheet1.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Add the new worksheet right after Sheet1 and name it "Result".
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Result"

# Header row.
$ws2.Range("A1").Value = "ISBN"
$ws2.Range("B1").Value = "Tittel"
$ws2.Range("C1").Value = "Forlag"

# Force column A (ISBN numbers) to be stored as text, like on Sheet1,
# so the long digit-strings aren't mangled into numbers.
$ws2.Range("A2:A5").NumberFormat = "@"

$ws2.Range("A2").Value = "9788203265082"
$ws2.Range("B2").Value = "Blokka"
$ws2.Range("C2").Value = "Aschehoug"

$ws2.Range("A3").Value = "9788205598126"
$ws2.Range("B3").Value = "Vinter i Applemore"
$ws2.Range("C3").Value = "Gyldendal"

$ws2.Range("A4").Value = "9788242184702"
$ws2.Range("B4").Value = "Alle utlendinger har lukka gardiner"
$ws2.Range("C4").Value = "Lydbokforlaget"

$ws2.Range("A5").Value = "9788249527496"
$ws2.Range("B5").Value = "Jeg plystrer i den mørke vinden"
$ws2.Range("C5").Value = "Forlaget Oktober"

# Turn the range into a real Excel Table ("ListObject").
$tbl = $ws2.ListObjects.Add(1, $ws2.Range("A1:C5"), $null, 1)
$tbl.Name = "MyTable"
$tbl.TableStyle = "TableStyleMedium2"
$tbl.ShowTableStyleRowStripes = $false

$wb.Save()
